## Adds the "1803" (18 March) minutes sheet, cloned from "1703" (17 March)
## with updated content, and tweaks the "1703" sheet view now that it is
## no longer the active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new sheet "1803" at the end of the tab strip, cloning the
#    previous week's sheet ("1703") as a starting point (values + styles).
# ---------------------------------------------------------------------
$src = $wb.Worksheets.Item("1703")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "1803"

$src.Range("A1:C21").Copy($ws.Range("A1:C21"))

# ---------------------------------------------------------------------
# 2. Content updates specific to the 18 March minutes.
# ---------------------------------------------------------------------

# Meeting date: 17th -> 18th March 2021.
$ws.Range("B1").Value = 44273

# "Complete stakeholder analysis" action is now finished.
$ws.Range("C10").Value = "Completed"

# Drop the second (now stale) blank action-review row, replacing it with
# new rows for the Docker-workshop action.
$ws.Range("A14").Value = "Do the Docker workshop"
$ws.Range("B14").Value = "All"
$ws.Range("C14").Value = "This week or next"

# New decisions made in this meeting.
$ws.Range("A17").Value = "Gone over scoping of project (see document Scoping in Prototyping folder)"
$ws.Range("A18").Value = "Feedback meeting booked for 31st March at 3.30pm"
$ws.Range("A19").Value = "Technical focus more on front-end side than back-end"

# ---------------------------------------------------------------------
# 3. Row heights / column widths for the new sheet.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 25.33
$ws.Columns.Item(2).ColumnWidth = 15
$ws.Columns.Item(3).ColumnWidth = 17.67

$ws.Rows.Item(7).RowHeight = 30
$ws.Rows.Item(8).RowHeight = 98.25
$ws.Rows.Item(9).RowHeight = 50.25
$ws.Rows.Item(10).RowHeight = 45
$ws.Rows.Item(17).RowHeight = 77.25
$ws.Rows.Item(18).RowHeight = 47.25
$ws.Rows.Item(19).RowHeight = 66.75

# ---------------------------------------------------------------------
# 4. View state: "1703" no longer the selected tab, scrolled down a bit
#    with the whole table selected; "1803" becomes the active tab with
#    C14 selected.
# ---------------------------------------------------------------------
$src.Range("A1:C21").Select()

$ws.Activate()
$ws.Range("C14").Select()

$wb.Save()
